$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two farm coordinate text labels (shared strings) in R25/S25
$ws.Range("R25").Value = "(38.18033471302274, 41.49317846476236)"
$ws.Range("S25").Value = "(37.65974498060651, 42.013768197178585)"

# Update Monte Carlo simulation result numbers for each load point row (2-24) and the TOTAL row (25)
$ws.Cells.Item(2, 6).Value = 7.277830834687032
$ws.Cells.Item(2, 7).Value = 2937
$ws.Cells.Item(2, 8).Value = 5.602715531912625
$ws.Cells.Item(2, 9).Value = 1.298982750995135
$ws.Cells.Item(2, 11).Value = 190.9504643962848
$ws.Cells.Item(2, 12).Value = 1069.841132698994
$ws.Cells.Item(2, 13).Value = 823.5991831911558
$ws.Cells.Item(3, 6).Value = 7.338183018174249
$ws.Cells.Item(3, 7).Value = 2938
$ws.Cells.Item(3, 8).Value = 5.647253847546622
$ws.Cells.Item(3, 9).Value = 1.299425033171163
$ws.Cells.Item(3, 11).Value = 163.7275541795666
$ws.Cells.Item(3, 12).Value = 924.6110602899553
$ws.Cells.Item(3, 13).Value = 711.5539847908744
$ws.Cells.Item(4, 6).Value = 7.345777214604732
$ws.Cells.Item(4, 7).Value = 2943
$ws.Cells.Item(4, 8).Value = 5.643493809793171
$ws.Cells.Item(4, 9).Value = 1.301636444051305
$ws.Cells.Item(4, 11).Value = 1.301636444051305
$ws.Cells.Item(4, 12).Value = 7.345777214604732
$ws.Cells.Item(4, 13).Value = 5.643493809793171
$ws.Cells.Item(5, 6).Value = 7.339855710347033
$ws.Cells.Item(5, 7).Value = 2948
$ws.Cells.Item(5, 8).Value = 5.629380515975115
$ws.Cells.Item(5, 9).Value = 1.303847854931446
$ws.Cells.Item(5, 11).Value = 1.303847854931446
$ws.Cells.Item(5, 12).Value = 7.339855710347033
$ws.Cells.Item(5, 13).Value = 5.629380515975115
$ws.Cells.Item(6, 6).Value = 8.146739711310651
$ws.Cells.Item(6, 7).Value = 2934
$ws.Cells.Item(6, 8).Value = 6.278043110863457
$ws.Cells.Item(6, 9).Value = 1.29765590446705
$ws.Cells.Item(6, 11).Value = 171.2905793896506
$ws.Cells.Item(6, 12).Value = 1075.369641893006
$ws.Cells.Item(6, 13).Value = 828.7016906339763
$ws.Cells.Item(7, 6).Value = 8.392871508974606
$ws.Cells.Item(7, 7).Value = 2991
$ws.Cells.Item(7, 8).Value = 6.344460876560208
$ws.Cells.Item(7, 9).Value = 1.322865988500664
$ws.Cells.Item(7, 11).Value = 194.4613003095975
$ws.Cells.Item(7, 12).Value = 1233.752111819267
$ws.Cells.Item(7, 13).Value = 932.6357488543506
$ws.Cells.Item(8, 6).Value = 8.449436424506292
$ws.Cells.Item(8, 7).Value = 3007
$ws.Cells.Item(8, 8).Value = 6.353234371735526
$ws.Cells.Item(8, 9).Value = 1.329942503317116
$ws.Cells.Item(8, 11).Value = 1.329942503317116
$ws.Cells.Item(8, 12).Value = 8.449436424506292
$ws.Cells.Item(8, 13).Value = 6.353234371735526
$ws.Cells.Item(9, 6).Value = 8.308972040212121
$ws.Cells.Item(9, 7).Value = 2946
$ws.Cells.Item(9, 8).Value = 6.376980917487985
$ws.Cells.Item(9, 9).Value = 1.30296329057939
$ws.Cells.Item(9, 11).Value = 102.9340999557718
$ws.Cells.Item(9, 12).Value = 656.4087911767575
$ws.Cells.Item(9, 13).Value = 503.7814924815508
$ws.Cells.Item(10, 6).Value = 10.33986683513122
$ws.Cells.Item(10, 7).Value = 2997
$ws.Cells.Item(10, 8).Value = 7.800613584995561
$ws.Cells.Item(10, 9).Value = 1.325519681556833
$ws.Cells.Item(10, 11).Value = 1.325519681556833
$ws.Cells.Item(10, 12).Value = 10.33986683513122
$ws.Cells.Item(10, 13).Value = 7.800613584995561
$ws.Cells.Item(11, 6).Value = 10.17095399398047
$ws.Cells.Item(11, 7).Value = 2936
$ws.Cells.Item(11, 8).Value = 7.832604557353489
$ws.Cells.Item(11, 9).Value = 1.298540468819107
$ws.Cells.Item(11, 11).Value = 98.68907563025211
$ws.Cells.Item(11, 12).Value = 772.9925035425159
$ws.Cells.Item(11, 13).Value = 595.2779463588652
$ws.Cells.Item(12, 6).Value = 10.38242590916942
$ws.Cells.Item(12, 7).Value = 3253
$ws.Cells.Item(12, 8).Value = 7.216312628537371
$ws.Cells.Item(12, 9).Value = 1.43874391862008
$ws.Cells.Item(12, 11).Value = 113.6607695709863
$ws.Cells.Item(12, 12).Value = 820.2116468243845
$ws.Cells.Item(12, 13).Value = 570.0886976544523
$ws.Cells.Item(13, 6).Value = 10.35474506092288
$ws.Cells.Item(13, 7).Value = 3243
$ws.Cells.Item(13, 8).Value = 7.219265674605809
$ws.Cells.Item(13, 9).Value = 1.434321096859797
$ws.Cells.Item(13, 11).Value = 109.0084033613445
$ws.Cells.Item(13, 12).Value = 786.960624630139
$ws.Cells.Item(13, 13).Value = 548.6641912700414
$ws.Cells.Item(14, 6).Value = 11.35741011934108
$ws.Cells.Item(14, 7).Value = 3533
$ws.Cells.Item(14, 8).Value = 7.268356716623321
$ws.Cells.Item(14, 9).Value = 1.562582927908005
$ws.Cells.Item(14, 11).Value = 1.562582927908005
$ws.Cells.Item(14, 12).Value = 11.35741011934108
$ws.Cells.Item(14, 13).Value = 7.268356716623321
$ws.Cells.Item(15, 6).Value = 10.3061842301787
$ws.Cells.Item(15, 7).Value = 3651
$ws.Cells.Item(15, 8).Value = 6.382438385218858
$ws.Cells.Item(15, 9).Value = 1.614772224679345
$ws.Cells.Item(15, 11).Value = 127.5670057496683
$ws.Cells.Item(15, 12).Value = 814.1885541841176
$ws.Cells.Item(15, 13).Value = 504.2126324322898
$ws.Cells.Item(16, 6).Value = 10.49999450621463
$ws.Cells.Item(16, 7).Value = 3713
$ws.Cells.Item(16, 8).Value = 6.393882999879149
$ws.Cells.Item(16, 9).Value = 1.6421937195931
$ws.Cells.Item(16, 11).Value = 1.6421937195931
$ws.Cells.Item(16, 12).Value = 10.49999450621463
$ws.Cells.Item(16, 13).Value = 6.393882999879149
$ws.Cells.Item(17, 6).Value = 10.36260656644125
$ws.Cells.Item(17, 7).Value = 3653
$ws.Cells.Item(17, 8).Value = 6.413866259710831
$ws.Cells.Item(17, 9).Value = 1.615656789031402
$ws.Cells.Item(17, 11).Value = 122.7899159663865
$ws.Cells.Item(17, 12).Value = 787.5580990495349
$ws.Cells.Item(17, 13).Value = 487.4538357380231
$ws.Cells.Item(18, 6).Value = 10.28541348685198
$ws.Cells.Item(18, 7).Value = 3647
$ws.Cells.Item(18, 8).Value = 6.37656152831706
$ws.Cells.Item(18, 9).Value = 1.613003095975232
$ws.Cells.Item(18, 11).Value = 1.613003095975232
$ws.Cells.Item(18, 12).Value = 10.28541348685198
$ws.Cells.Item(18, 13).Value = 6.37656152831706
$ws.Cells.Item(19, 6).Value = 11.38345753309828
$ws.Cells.Item(19, 7).Value = 3987
$ws.Cells.Item(19, 8).Value = 6.455479679542314
$ws.Cells.Item(19, 9).Value = 1.763379035824856
$ws.Cells.Item(19, 11).Value = 1.763379035824856
$ws.Cells.Item(19, 12).Value = 11.38345753309828
$ws.Cells.Item(19, 13).Value = 6.455479679542314
$ws.Cells.Item(20, 6).Value = 11.97047365771821
$ws.Cells.Item(20, 7).Value = 3715
$ws.Cells.Item(20, 8).Value = 7.285394600296329
$ws.Cells.Item(20, 9).Value = 1.643078283945157
$ws.Cells.Item(20, 11).Value = 129.8031844316674
$ws.Cells.Item(20, 12).Value = 945.6674189597383
$ws.Cells.Item(20, 13).Value = 575.54617342341
$ws.Cells.Item(21, 6).Value = 12.14361176300139
$ws.Cells.Item(21, 7).Value = 3777
$ws.Cells.Item(21, 8).Value = 7.269448291275123
$ws.Cells.Item(21, 9).Value = 1.670499778858912
$ws.Cells.Item(21, 11).Value = 1.670499778858912
$ws.Cells.Item(21, 12).Value = 12.14361176300139
$ws.Cells.Item(21, 13).Value = 7.269448291275123
$ws.Cells.Item(22, 6).Value = 11.96337606491591
$ws.Cells.Item(22, 7).Value = 3707
$ws.Cells.Item(22, 8).Value = 7.296788044989174
$ws.Cells.Item(22, 9).Value = 1.639540026536931
$ws.Cells.Item(22, 11).Value = 1.639540026536931
$ws.Cells.Item(22, 12).Value = 11.96337606491591
$ws.Cells.Item(22, 13).Value = 7.296788044989174
$ws.Cells.Item(23, 6).Value = 11.95846416864866
$ws.Cells.Item(23, 7).Value = 3709
$ws.Cells.Item(23, 8).Value = 7.289859122489786
$ws.Cells.Item(23, 9).Value = 1.640424590888987
$ws.Cells.Item(23, 11).Value = 124.672268907563
$ws.Cells.Item(23, 12).Value = 908.843276817298
$ws.Cells.Item(23, 13).Value = 554.0292933092237
$ws.Cells.Item(24, 6).Value = 12.9492694175518
$ws.Cells.Item(24, 7).Value = 3998
$ws.Cells.Item(24, 8).Value = 7.32323615634933
$ws.Cells.Item(24, 9).Value = 1.768244139761168
$ws.Cells.Item(24, 11).Value = 1.768244139761168
$ws.Cells.Item(24, 12).Value = 12.9492694175518
$ws.Cells.Item(24, 13).Value = 7.32323615634933
$ws.Cells.Item(25, 11).Value = 1.408685554570629
$ws.Cells.Item(25, 12).Value = 9.222706957701837
$ws.Cells.Item(25, 13).Value = 6.547030263622559
$ws.Cells.Item(25, 14).Value = 39.83675658889262
$ws.Cells.Item(25, 15).Value = 2261
$ws.Cells.Item(25, 17).Value = 0.02121440785693962
